# Applies the edits described in the commit diff to the 'Pina' worksheet
# (and updates the active-cell selection).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pina")

# --- Simple numeric corrections in the top summary block (rows 2-17, 37, 46) ---
$ws.Range("E2").Value = 0
$ws.Range("E3").Value = 24668.7
$ws.Range("E4").Value = 0.35
$ws.Range("D5").Value = 814393.5
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E8").Value = 896.27
$ws.Range("D9").Value = 6279138.02
$ws.Range("E10").Value = 34976.21
$ws.Range("D11").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 2385537.27
$ws.Range("D16").Value = 1152391.92
$ws.Range("D17").Value = 406190.3
$ws.Range("D37").Value = 346621.65
$ws.Range("D46").Value = 0

# --- Rows 63-76: new supplier/OC concept text in column C + paid amount in column E ---
$ws.Range("C63").Value = 'EARTHCROP S.A OC 8433 (COMPRA DE MIRAGE )'
$ws.Range("E63").Value = 3450
$ws.Range("C64").Value = 'ENLACE AGROPECUARIO OC 8434 (COMPRA DE TECNOSILIX)'
$ws.Range("E64").Value = 384
$ws.Range("C65").Value = 'UPL COSTA RICA S.A OC 8438 (COMPRA DE BELTANOL)'
$ws.Range("E65").Value = 700
$ws.Range("C66").Value = 'David Salazar Serrano OC 8415 (COMORA DE FILTRO ACTIVO 03-003)'
$ws.Range("E66").Value = 59.27
$ws.Range("C67").Value = '3-101809221 OC 8421 (COMPRA DE ARRANCADOR ACTIVO 00-025)'
$ws.Range("E67").Value = 521.66
$ws.Range("C68").Value = 'David Salazar Serrano OC 8440 (COMPRA DE ELECTROVALVULA ACTIVO 03-001)'
$ws.Range("E68").Value = 419.98
$ws.Range("C69").Value = 'ECONOVA OC 8457 (COMPRA DE MATRICES PARA EL LABORATORIO)'
$ws.Range("E69").Value = 59.8499
$ws.Range("C70").Value = 'David Salazar Serrano OC 8461 (COMPRA DE SELLOS, ROLL Y MOTOR HIDRAULICO ACTIVO 03-003)'
$ws.Range("E70").Value = 1923.48
$ws.Range("C71").Value = 'ABONOS DEL PACIFICO S.A OC 8462 (COMPRA DE UREA, CALCIO, AMONIO, MAGNESIO)'
$ws.Range("E71").Value = 27837.425
$ws.Range("C72").Value = 'BIOCAMPO OC 8464 (COMPRA DE SILWET APLICACIONES SEMANA 07)'
$ws.Range("E72").Value = 960
$ws.Range("C73").Value = 'ENLACE AGROPECUARIO OC 8469 (COMPRA DE PROCTESOL SOLAR SEMANA 07)'
$ws.Range("E73").Value = 11150
$ws.Range("C74").Value = 'UPL COSTA RICA S.A OC 8472 (COMPEA DE BELTANOL APLICACIONES SEMANA 07)'
$ws.Range("E74").Value = 700
$ws.Range("C75").Value = 'YARA COSTA RICA S.R.L OC 8473 (COMPRA DE AZUTECK, COMPLEX Y WHITE BALANCE)'
$ws.Range("E75").Value = 29175
$ws.Range("C76").Value = 'ChemTica Internacional S.A OC 8494 (COMPRA DE FEROMONAS)'
$ws.Range("E76").Value = 310

# --- Rows 86-107: new supplier/OC concept text in column C + paid amount in column D ---
$ws.Range("C86").Value = 'KEVIN TIJERINO ARAGON OC 8369 (SERVICIO DE MANTENIMIENTO Y REPARACION DE AC LABORATORIOY CASA)'
$ws.Range("D86").Value = 370000
$ws.Range("C87").Value = 'DIXIE ARLEY JIMENEZ OC 8142 (COMPRA DE ROLL ACTIVO 07-006)'
$ws.Range("C88").Value = 'Alquimia Industrial S,A OC 8449 (COMPRA DE LIMPIEZA PARA LABORATORIO Y COCINA)'
$ws.Range("D88").Value = 34552.53
$ws.Range("C89").Value = 'GRUPO WEFEL OC 8122 (COMPRA DE REPUESTOS ACTIVO 00-019)'
$ws.Range("D89").Value = 353147.41
$ws.Range("C90").Value = 'BATERIAS LA BODEGUITA OC 8453 (COMPRA DE BATERIA ACTIVO 00-019,00-021 Y 00-041)'
$ws.Range("D90").Value = 199860
$ws.Range("C91").Value = 'GRUPO WEFEL OC 8454 (COMPRA DE REPUESTOS ACTIVO 00-018)'
$ws.Range("D91").Value = 871065.64
$ws.Range("C92").Value = 'LUIS CARLOS CASTILLO CASTILLO OC 8459 (COMPRA DE MECATE PARA SEMILLA Y SIEMBRA)'
$ws.Range("D92").Value = 26813.48
$ws.Range("C93").Value = 'CENTRAL DE BOLSAS PLASTICAS OC 8460 (COMPRA DE BOLSAS PARA EL LABORATORIO)'
$ws.Range("D93").Value = 108265.33
$ws.Range("C94").Value = 'Green Go S,A OC 8476 (COMPRA DE SET DE EMPAQUES CARGADOR MILITAR )'
$ws.Range("D94").Value = 157831.85
$ws.Range("C95").Value = 'GRUPO WEFEL OC 8480 (COMPRA DE FILTRO PARA MASEY FERGUSON)'
$ws.Range("D95").Value = 37272.6
$ws.Range("C96").Value = 'Cristian Arroyo Zuñiga OC 8481 (COMPRA DE MANGAS ARMI)'
$ws.Range("D96").Value = 32000
$ws.Range("C97").Value = 'GRUPO WEFEL OC 8484 (COMPRA DE ACOPLES PARA TRACTORES)'
$ws.Range("D97").Value = 189412.36
$ws.Range("C98").Value = 'CENTRAL DE BOLSAS PLASTICAS OC 8487 (COMPRA DE BOLSA NEGRA PARA CONTROL DE CABEZA ROJA)'
$ws.Range("D98").Value = 26175.05
$ws.Range("C99").Value = 'GASOLINERA KATIRA S.A. OC 8489 (COMPRA DE GASOLINA BODEGA)'
$ws.Range("D99").Value = 177300
$ws.Range("C100").Value = 'Liga Agrícola Industrial de la Caña de Azúcar OC 8490 (COMPRA DE MIEL )'
$ws.Range("D100").Value = 2391680
$ws.Range("C101").Value = 'LA CASA DE LOS ROLES S.A. OC 8491 (COMPRA DE ROLES Y RETENEDORES CARRETAS)'
$ws.Range("D101").Value = 176320
$ws.Range("C102").Value = 'Llantas Importadas Llisa S,A OC 8492 (COMPRA DE LLANTAS Y NEUMATICOS 12,5L15 CARRETAS)'
$ws.Range("D102").Value = 342800
$ws.Range("C103").Value = 'LLANTAS DEL PACIFICO SAN CARLOS S.A OC 8493 (COMPRA DE NEUMATICOS PARA LA BODEGA )'
$ws.Range("D103").Value = 124480.02
$ws.Range("C104").Value = 'Alquimia Industrial S,A OC 8496 (COMPRA DE TOALLS, DTERGENTE PAPEL Y ESCOBA PARA TALLER)'
$ws.Range("D104").Value = 31901.0069
$ws.Range("C105").Value = 'Alquimia Industrial S,A OC 8495 (COMPRA DE PALAS Y ESCOBA BODEGA INSUMOS)'
$ws.Range("D105").Value = 8655.01
$ws.Range("C106").Value = 'Servicios Electrónicos Azocar Ltda. OC 8497 (COMPRA DE ROMANA PARA BODEGA)'
$ws.Range("D106").Value = 222893.8
$ws.Range("C107").Value = 'Tecnología Extrema San Carlos Limitada OC 8498 (COMPRA DE 2 COMPUTADORAS PARA USO PREPARACION Y SEMILLA YS IEMBRA)'
$ws.Range("D107").Value = 514000

# --- Misc single-cell numeric corrections further down the sheet ---
$ws.Range("E221").Value = 283
$ws.Range("E224").Value = 673.41
$ws.Range("E226").Value = 150

# --- Update the active-cell selection to match the saved worksheet view ---
$ws.Range("F10").Select()
